# Auto-generated edit script: updates cryptos Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = "29.977.02"
$ws.Cells.Item(2, 5).Value2 = "  -0.20%  "

$ws.Cells.Item(3, 4).Value2 = "1.896.56"
$ws.Cells.Item(3, 5).Value2 = "  -0.80%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value2 = "0.8433"
$ws.Cells.Item(5, 4).ClearFormats()

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value2 = "241.65"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value2 = "  -0.31%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value2 = "0.9998"
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value2 = "  +0.03%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value2 = "0.3309"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value2 = "  +3.88%  "

$ws.Cells.Item(9, 5).Value2 = "  +1.53%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value2 = "0.07050"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value2 = "  +1.42%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value2 = "0.08085"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value2 = "  +0.57%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value2 = "0.7592"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value2 = "  +1.11%  "

$ws.Cells.Item(13, 4).Value2 = "1.899.75"
$ws.Cells.Item(13, 5).Value2 = "  -0.57%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value2 = "5.264"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value2 = "  +0.54%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value2 = "92.26"
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).Value2 = "  -1.49%  "

$ws.Cells.Item(16, 4).Value2 = "29.978.16"
$ws.Cells.Item(16, 5).Value2 = "  -0.20%  "

$ws.Cells.Item(17, 5).Value2 = "  +0.55%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value2 = "5.879"
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).Value2 = "  -1.90%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value2 = "244.64"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value2 = "  -2.38%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value2 = "0.000007774"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value2 = "  -0.41%  "

$ws.Cells.Item(21, 5).Value2 = "  -0.11%  "

$ws.Cells.Item(22, 4).Value2 = "2.151.39"
$ws.Cells.Item(22, 5).Value2 = "  -0.48%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value2 = "0.9998"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value2 = "  -0.06%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value2 = "6.990"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value2 = "  +0.00%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value2 = "0.1747"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value2 = "  +28.28%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value2 = "9.261"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value2 = "  -0.69%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value2 = "166.13"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value2 = "  -1.59%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value2 = "18.88"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value2 = "  -0.82%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value2 = "2.110"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value2 = "  +1.91%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value2 = "1.362"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value2 = "  -2.04%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value2 = "1.518"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value2 = "  -0.40%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value2 = "0.05824"
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value2 = "  +7.81%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value2 = "4.296"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value2 = "  -1.66%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value2 = "4.078"
$ws.Cells.Item(34, 4).ClearFormats()

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value2 = "1.273"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value2 = "  +0.69%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value2 = "0.7319"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value2 = "  -1.36%  "

$ws.Cells.Item(37, 5).Value2 = "  -0.16%  "

$ws.Cells.Item(38, 5).Value2 = "  -0.72%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value2 = "2.774"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value2 = "  -0.53%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value2 = "0.4443"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value2 = "  -0.82%  "

$ws.Cells.Item(41, 5).Value2 = "  -0.92%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value2 = "5.881"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value2 = "  -4.90%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value2 = "0.8424"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value2 = "  +1.16%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value2 = "0.9993"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value2 = "  -0.10%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value2 = "1.890"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value2 = "  -1.37%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value2 = "101.64"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value2 = "  +1.01%  "

$ws.Cells.Item(47, 4).Value2 = "1.010.33"
$ws.Cells.Item(47, 5).Value2 = "  +4.80%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value2 = "7.580"
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value2 = "  -1.19%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value2 = "9.815"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value2 = "  -0.59%  "

$ws.Cells.Item(50, 4).Value2 = "2.047.39"
$ws.Cells.Item(50, 5).Value2 = "  -0.54%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value2 = "35.94"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value2 = "  -1.42%  "
